$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple per-cell updates for price/volume columns (rows with no coin identity change) ---
# Cells whose new Price text looks like a plain number are forced to a Text number format first,
# so Excel keeps storing the literal text (matching the original inlineStr string cells) instead of
# silently converting it to a numeric value (which would drop formatting like trailing zeros).
$ws.Range("D2").Value = "79.699.71"
$ws.Range("E2").Value = "  +5.07%  "
$ws.Range("D3").Value = "3.208.31"
$ws.Range("E3").Value = "  +6.86%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.88"
$ws.Range("E5").Value = "  +7.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "638.77"
$ws.Range("E6").Value = "  +4.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.264"
$ws.Range("E7").Value = "  +30.43%  "
$ws.Range("E9").Value = "  +10.58%  "
$ws.Range("D10").Value = "3.207.78"
$ws.Range("E10").Value = "  +6.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.613"
$ws.Range("E11").Value = "  +40.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000269"
$ws.Range("E12").Value = "  +41.98%  "
$ws.Range("E13").Value = "  +3.76%  "
$ws.Range("E14").Value = "  +5.27%  "
$ws.Range("D15").Value = "3.799.37"
$ws.Range("E15").Value = "  +6.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.70"
$ws.Range("E16").Value = "  +14.19%  "
$ws.Range("D17").Value = "79.480.97"
$ws.Range("E17").Value = "  +4.66%  "
$ws.Range("D18").Value = "3.205.20"
$ws.Range("E18").Value = "  +6.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.64"
$ws.Range("E19").Value = "  +9.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.43"
$ws.Range("E20").Value = "  +6.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.01"
$ws.Range("E21").Value = "  +29.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "447.31"
$ws.Range("E22").Value = "  +19.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.32"
$ws.Range("E23").Value = "  +22.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.88"
$ws.Range("E24").Value = "  +14.20%  "
$ws.Range("D25").Value = "3.366.33"
$ws.Range("E25").Value = "  +6.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "77.83"
$ws.Range("E26").Value = "  +7.84%  "
$ws.Range("E27").Value = "  +13.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("E29").Value = "  +18.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.29"
$ws.Range("E30").Value = "  +13.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "564.99"
$ws.Range("E32").Value = "  +16.00%  "
$ws.Range("E33").Value = "  +11.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.158"
$ws.Range("E34").Value = "  +33.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.05"
$ws.Range("E35").Value = "  +7.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.29"
$ws.Range("E36").Value = "  +14.43%  "
$ws.Range("E37").Value = "  +21.67%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.415"
$ws.Range("E39").Value = "  +10.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.59"
$ws.Range("E40").Value = "  +1.14%  "

# --- Rows 41 and 42 swapped coin identities (RenderToken <-> WhiteBITCoin), with updated price/volume ---
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.76"
$ws.Range("E41").Value = "  +14.24%  "

$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.30"
$ws.Range("E42").Value = "  +1.40%  "

# --- Remaining per-cell updates for price/volume columns ---
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "192.19"
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  +13.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.74"
$ws.Range("E46").Value = "  +14.57%  "
$ws.Range("E47").Value = "  +5.36%  "
$ws.Range("E48").Value = "  +9.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "43.20"
$ws.Range("E49").Value = "  +5.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.32"
$ws.Range("E50").Value = "  +12.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.95"
$ws.Range("E51").Value = "  +17.98%  "
